$d = $word.ActiveDocument

# "Testing:" section - the empty bulleted list item right after it
# gets the new "Functionality of game works" text.
$testingBullet = $d.Paragraphs.Item(7)
$testingBullet.Range.Text = "Functionality of game works"

# "Bug Reports:" section - the empty bulleted list item right after it
# gets the new "Winning a game still counts towards your losses" text.
$bugBullet = $d.Paragraphs.Item(9)
$bugBullet.Range.Text = "Winning a game still counts towards your losses"

# Add a brand-new bulleted list item right after that one for the
# second bug report.
$bugBullet.Range.InsertParagraphAfter()
$newBullet = $d.Paragraphs.Item(10)
$newBullet.Range.Text = "See all bots button is not working"
